# Update "想去人数" (interest count) values in column F for the "展览"
# and "全部类型" worksheets, as published by the gh-pages data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    2  = 196
    3  = 5422
    6  = 26
    7  = 623
    9  = 1061
    11 = 1494
    12 = 4538
    14 = 201
    15 = 173
    16 = 99
    17 = 3543
    18 = 179
    19 = 1117
    22 = 204
    23 = 25
    24 = 137
    25 = 48
    28 = 319
    33 = 34
}
foreach ($row in $updates1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates1[$row]
}

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    2  = 196
    4  = 5422
    6  = 55
    7  = 26
    8  = 623
    12 = 1494
    13 = 4538
    15 = 201
    16 = 173
    17 = 99
    18 = 3543
    19 = 179
    20 = 1117
    23 = 204
    24 = 25
    25 = 137
    26 = 48
    28 = 73
    29 = 319
    34 = 34
}
foreach ($row in $updates4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updates4[$row]
}
